$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "urn:eu.toop.process.twophasedrequestresponse `n"
$ws.Range("A4").Value = "TOOP Two Phased Request Response for Documents"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = $false
